$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of execution parameters (row 5) - an ID for Tareq
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "synthetic"
$ws.Range("C5").Value = "original"
$ws.Range("D5").Value = $false
$ws.Range("F5").Value = ".vtp"
$ws.Range("G5").Value = 200
$ws.Range("H5").Value = "M"
$ws.Range("I5").Value = $false
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "manually"
$ws.Range("L5").Value = '{"TRAGION_RIGHT": "4869", "TRAGION_LEFT": "2431", "NASION": "9396"}'
$ws.Range("M5").Value = $true
$ws.Range("N5").Value = $true
$ws.Range("O5").Value = $true

# Move active selection to K5, matching the post-edit cursor position
$ws.Range("K5").Select()
